$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing row (111) down into the two new rows
# so the new cells inherit the same styles (bordered/centered index col,
# date-formatted match-date col) without creating duplicate style entries.
$ws.Range("A111:V111").Copy($ws.Range("A112:V112"))
$ws.Range("A111:V111").Copy($ws.Range("A113:V113"))

# New row 112 (Indice 111)
$ws.Range("A112").Value = 111
$ws.Range("B112").Value = "costa-rica"
$ws.Range("C112").Value = "primera-division"
$ws.Range("D112").Value = "2023-2024"
$ws.Range("E112").Value = 45239.13541666666
$ws.Range("F112").Value = "Puntarenas FC"
$ws.Range("G112").Value = 4
$ws.Range("H112").Value = "Sporting San Jose"
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2.13
$ws.Range("K112").Value = "05/11/2023 00:12"
$ws.Range("L112").Value = 2.18
$ws.Range("M112").Value = "09/11/2023 03:11"
$ws.Range("N112").Value = 3.37
$ws.Range("O112").Value = "05/11/2023 00:12"
$ws.Range("P112").Value = 3.48
$ws.Range("Q112").Value = "09/11/2023 03:11"
$ws.Range("R112").Value = 3.45
$ws.Range("S112").Value = "05/11/2023 00:12"
$ws.Range("T112").Value = 3.35
$ws.Range("U112").Value = "09/11/2023 03:11"
$ws.Range("V112").Value = "https://www.betexplorer.com/football/costa-rica/primera-division/puntarenas-fc-sporting-san-jose/UcDN9TAH/"

# New row 113 (Indice 112)
$ws.Range("A113").Value = 112
$ws.Range("B113").Value = "costa-rica"
$ws.Range("C113").Value = "primera-division"
$ws.Range("D113").Value = "2023-2024"
$ws.Range("E113").Value = 45239.14583333334
$ws.Range("F113").Value = "Herediano"
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = "Cartagines"
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1.45
$ws.Range("K113").Value = "05/11/2023 22:12"
$ws.Range("L113").Value = 1.65
$ws.Range("M113").Value = "09/11/2023 03:21"
$ws.Range("N113").Value = 4.41
$ws.Range("O113").Value = "05/11/2023 22:12"
$ws.Range("P113").Value = 4.16
$ws.Range("Q113").Value = "09/11/2023 03:21"
$ws.Range("R113").Value = 5.98
$ws.Range("S113").Value = "05/11/2023 22:12"
$ws.Range("T113").Value = 4.92
$ws.Range("U113").Value = "09/11/2023 03:21"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/costa-rica/primera-division/herediano-cartagines/tO17DBBh/"
